$wb = $excel.ActiveWorkbook

# The workbook currently has two sheets, in this order:
#   1) "Step-by-step statistics"
#   2) "Preliminary statistics"
# We need to insert two brand-new sheets *before* the existing ones, so the
# final left-to-right tab order becomes:
#   1) "DAU statistics"             (new)
#   2) "MAU statistics"             (new)
#   3) "Step-by-step statistics"    (untouched)
#   4) "Preliminary statistics"     (untouched)

$firstExisting = $wb.Worksheets.Item(1)

# NOTE: a worksheet handle returned by Worksheets.Add() tracks a *position*,
# not a stable object identity - inserting another sheet before it re-points
# the existing variable to whatever now occupies that slot. So each new
# sheet must be fully populated (named + filled in) before the next Add()
# call runs.

# --- MAU statistics sheet (inserted first, ends up second) ---------------
$mau = $wb.Worksheets.Add($firstExisting)
$mau.Name = "MAU statistics"

$mau.Range("A1").Value = "Month"
$mau.Range("B1").Value = "MAU"

# Force A2 to be stored as literal text ("01.01.2018") instead of letting
# Excel auto-convert the date-looking string into a date serial number:
# flip the cell to Text format, assign the value, then restore the cell's
# style to Normal so no stray number formatting lingers on it.
$mau.Range("A2").NumberFormat = "@"
$mau.Range("A2").Value = "01.01.2018"
$mau.Range("A2").Style = "Normal"

$mau.Range("B2").Value = 14831

# --- DAU statistics sheet (inserted second, ends up first) ---------------
$dau = $wb.Worksheets.Add($firstExisting)
$dau.Name = "DAU statistics"

$dau.Range("A1").Value = "Date"
$dau.Range("B1").Value = "DAU"

$dau.Range("A2").NumberFormat = "@"
$dau.Range("A2").Value = "01.01.2018"
$dau.Range("A2").Style = "Normal"

$dau.Range("B2").Value = 14831
